# Add a new "2022" column (K) to the poverty-rate table, mirroring the
# formatting of the existing "2020" column (I) and filling in the new
# year's figures for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Clone the formatting of column I (2020) into the new column K so the
#    new column looks consistent with the rest of the table (same number
#    format, font, fills, borders, alignment) for header row 4 through the
#    last data row 46.
$ws.Range("I4:I46").Copy()
$ws.Range("K4:K46").PasteSpecial(-4122)

# 2) Header for the new year.
$ws.Cells.Item(4, 11).Value = 2022

# 3) New 2022 values for every row that carries data (rows that are purely
#    section/spacer headers - 6, 9, 22, 25 - stay empty, matching column I).
$values = @{
    5  = 24.2
    7  = 25.5
    8  = 22.3
    10 = 18
    11 = 18.9
    12 = 21.3
    13 = 30.2
    14 = 31.7
    15 = 34.1
    16 = 25.8
    17 = 20
    18 = 12.1
    19 = 10.3
    20 = 15.1
    21 = 12.1
    23 = 25.9
    24 = 23.2
    26 = 25.9
    27 = 48.3
    28 = 24.3
    29 = 28.1
    30 = 25.8
    31 = 27.1
    32 = 20.7
    33 = 24.3
    34 = 19.4
    35 = 7.5
    36 = 11.4
    37 = 36.5
    38 = 17.8
    39 = 20.3
    40 = 20.5
    41 = 32.2
    42 = 23.2
    43 = 23.8
    44 = 21
    45 = 18
    46 = 3.2
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 11).Value = $values[$row]
}

# 4) Move the active selection the way the author's workbook ended up
#    (selection moved down while entering the new column's data).
$ws.Range("L12").Select()
